# Update the "Inscricoes" worksheet with revised enrollment counts.
# Targets columns E (Inscritos), F (Pagos), and H (Inscricoes homologadas)
# for the specific rows whose figures changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E4").Value = 43
$ws.Range("E5").Value = 133
$ws.Range("F5").Value = 90
$ws.Range("H5").Value = 90
$ws.Range("E10").Value = 502
$ws.Range("F10").Value = 250
$ws.Range("H10").Value = 250
$ws.Range("E11").Value = 332
$ws.Range("E12").Value = 492
$ws.Range("F12").Value = 262
$ws.Range("H12").Value = 262
$ws.Range("E14").Value = 123
$ws.Range("E16").Value = 196
$ws.Range("F16").Value = 97
$ws.Range("H16").Value = 97
$ws.Range("E17").Value = 95
$ws.Range("F17").Value = 50
$ws.Range("H17").Value = 50
$ws.Range("E22").Value = 162
$ws.Range("E23").Value = 196
$ws.Range("F23").Value = 92
$ws.Range("H23").Value = 92
$ws.Range("E24").Value = 206
$ws.Range("E25").Value = 255
$ws.Range("E26").Value = 152
$ws.Range("F26").Value = 93
$ws.Range("H26").Value = 93
$ws.Range("E28").Value = 191
$ws.Range("F28").Value = 73
$ws.Range("H28").Value = 73
$ws.Range("E33").Value = 281
$ws.Range("E34").Value = 210
$ws.Range("F34").Value = 139
$ws.Range("H34").Value = 139
$ws.Range("E35").Value = 144
$ws.Range("E37").Value = 151
$ws.Range("E41").Value = 380
$ws.Range("E42").Value = 364
$ws.Range("F42").Value = 201
$ws.Range("H42").Value = 201
$ws.Range("E43").Value = 112
$ws.Range("E44").Value = 303
$ws.Range("F44").Value = 154
$ws.Range("H44").Value = 154
$ws.Range("E45").Value = 136
$ws.Range("E47").Value = 435
$ws.Range("F47").Value = 215
$ws.Range("H47").Value = 215
$ws.Range("E48").Value = 197
$ws.Range("E49").Value = 277
$ws.Range("E50").Value = 235
$ws.Range("E51").Value = 224
